$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") date serial bumps from 45184 to 45186 for rows 2..27.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# 2) For the first four data rows (2..5), the HYPERLINK() formulas in columns
#    S, T, V, W, X, Y gain a second argument - the friendly link text, which
#    equals the designation in column A of that row.
$folders = @{
    "S" = @{ Path = "artfynd";        Ext = ".xlsx" }
    "T" = @{ Path = "kartor";         Ext = ".png"  }
    "V" = @{ Path = "klagomål";       Ext = ".docx" }
    "W" = @{ Path = "klagomålsmail";  Ext = ".docx" }
    "X" = @{ Path = "tillsyn";        Ext = ".docx" }
    "Y" = @{ Path = "tillsynsmail";   Ext = ".docx" }
}

for ($r = 2; $r -le 5; $r++) {
    $designation = $ws.Range("A" + $r).Text

    foreach ($col in "S", "T", "V", "W", "X", "Y") {
        $info = $folders[$col]
        $url = "https://klasma.github.io/Logging_YSTAD/" + $info.Path + "/" + $designation + $info.Ext
        $formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
        $ws.Range($col + $r).Formula = $formula
    }
}
